$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18
$ws.Range("B18").Value = "-"
$ws.Range("C18").Value = "-"
$ws.Range("D18").Value = "-"
$ws.Range("E18").Value = "['MEC-1NA-Trat. Termicos', -, 'MEC-1NB-T.M. Metalicos', -]"
$ws.Range("F18").Value = "['ELM-1NA-Tecnologia dos Materiais.', 'ELM-1NA-Tecnologia dos Materiais.']"

# Row 19
$ws.Range("B19").Value = "-"
$ws.Range("C19").Value = "-"
$ws.Range("D19").Value = "-"
$ws.Range("E19").Value = "['MEC-1NA-Trat. Termicos', -, 'MEC-1NB-T.M. Metalicos', -]"
$ws.Range("F19").Value = "['ELM-1NA-Tecnologia dos Materiais.', 'ELM-1NA-Tecnologia dos Materiais.']"

# Row 20
$ws.Range("B20").Value = "-"
$ws.Range("C20").Value = "-"
$ws.Range("E20").Value = "['MEC-1NA-Trat. Termicos', -, 'MEC-1NB-T.M. Metalicos', -]"
$ws.Range("F20").Value = "-"

# Row 21
$ws.Range("B21").Value = "-"
$ws.Range("E21").Value = "['MEC-1NA-Trat. Termicos', -, 'MEC-1NB-T.M. Metalicos', -]"
